$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Schedule" (sheet1.xml): update row 3 and append new row 4
# ------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Schedule")

$fmtA_sched = $wsSchedule.Cells.Item(2, 1).NumberFormat
$fmtB_sched = $wsSchedule.Cells.Item(2, 2).NumberFormat

$scheduleData = @(
  @(2, 46073.04166666666, 46073.20833333334, 4, 15.12, 636.716379, 42.11087162698413),
  @(3, 46073.375, 46073.66666666666, 7, 26.46, 578.43281925, 21.86065076530612),
  @(4, 46074.29166666666, 46074.75, 11, 41.58, 443.86831125, 10.67504356060606)
)

foreach ($row in $scheduleData) {
    $r = $row[0]
    $wsSchedule.Cells.Item($r, 1).Value = $row[1]
    $wsSchedule.Cells.Item($r, 1).NumberFormat = $fmtA_sched
    $wsSchedule.Cells.Item($r, 2).Value = $row[2]
    $wsSchedule.Cells.Item($r, 2).NumberFormat = $fmtB_sched
    $wsSchedule.Cells.Item($r, 3).Value = $row[3]
    $wsSchedule.Cells.Item($r, 4).Value = $row[4]
    $wsSchedule.Cells.Item($r, 5).Value = $row[5]
    $wsSchedule.Cells.Item($r, 6).Value = $row[6]
}

# ------------------------------------------------------------------
# Sheet "Detailed" (sheet2.xml): update rows 2-49 and append new
# rows 50-97
# ------------------------------------------------------------------
$wsDetailed = $wb.Worksheets.Item("Detailed")

$fmtA_det = $wsDetailed.Cells.Item(2, 1).NumberFormat
$fmtD_det = $wsDetailed.Cells.Item(2, 4).NumberFormat

$detailedData = @(
  @(2, 46073, 84.79000000000001, "historical", 46073, "OFF"),
  @(3, 46073.02083333334, 89.06780999999999, "historical", 46073, "OFF"),
  @(4, 46073.04166666666, 84.79000000000001, "historical", 46073, "ON"),
  @(5, 46073.0625, 78, "historical", 46073, "ON"),
  @(6, 46073.08333333334, 83.25254, "historical", 46073, "ON"),
  @(7, 46073.10416666666, 85.65000000000001, "historical", 46073, "ON"),
  @(8, 46073.125, 78, "historical", 46073, "ON"),
  @(9, 46073.14583333334, 79.6999, "historical", 46073, "ON"),
  @(10, 46073.16666666666, 78, "historical", 46073, "ON"),
  @(11, 46073.1875, 85.65000000000001, "historical", 46073, "ON"),
  @(12, 46073.20833333334, 91.44624, "historical", 46073, "OFF"),
  @(13, 46073.22916666666, 110.13518, "historical", 46073, "OFF"),
  @(14, 46073.25, 120.66799, "historical", 46073, "OFF"),
  @(15, 46073.27083333334, 138.42, "historical", 46073, "OFF"),
  @(16, 46073.29166666666, 119.39764, "historical", 46073, "OFF"),
  @(17, 46073.3125, 105, "historical", 46073, "OFF"),
  @(18, 46073.33333333334, 79.95, "historical", 46073, "OFF"),
  @(19, 46073.35416666666, 69.5744, "historical", 46073, "OFF"),
  @(20, 46073.375, 53.43962, "historical", 46073, "ON"),
  @(21, 46073.39583333334, 51.45378, "historical", 46073, "ON"),
  @(22, 46073.41666666666, 36.63752, "historical", 46073, "ON"),
  @(23, 46073.4375, 36.06, "historical", 46073, "ON"),
  @(24, 46073.45833333334, 36.06, "historical", 46073, "ON"),
  @(25, 46073.47916666666, 36.06, "historical", 46073, "ON"),
  @(26, 46073.5, 36.06, "historical", 46073, "ON"),
  @(27, 46073.52083333334, 36.06, "historical", 46073, "ON"),
  @(28, 46073.54166666666, 36.06, "historical", 46073, "ON"),
  @(29, 46073.5625, 36.06, "historical", 46073, "ON"),
  @(30, 46073.58333333334, 36.06, "historical", 46073, "ON"),
  @(31, 46073.60416666666, 52.2928, "historical", 46073, "ON"),
  @(32, 46073.625, 53.90061, "historical", 46073, "ON"),
  @(33, 46073.64583333334, 57.0601, "historical", 46073, "ON"),
  @(34, 46073.66666666666, 57.31, "forecast", 46073, "OFF"),
  @(35, 46073.6875, 58.95134, "forecast", 46073, "OFF"),
  @(36, 46073.70833333334, 80.66354, "forecast", 46073, "OFF"),
  @(37, 46073.72916666666, 79.95, "forecast", 46073, "OFF"),
  @(38, 46073.75, 90.31362, "forecast", 46073, "OFF"),
  @(39, 46073.77083333334, 108.01, "forecast", 46073, "OFF"),
  @(40, 46073.79166666666, 114.5425, "forecast", 46073, "OFF"),
  @(41, 46073.8125, 109.78154, "forecast", 46073, "OFF"),
  @(42, 46073.83333333334, 120.01745, "forecast", 46073, "OFF"),
  @(43, 46073.85416666666, 108.89, "forecast", 46073, "OFF"),
  @(44, 46073.875, 108.01, "forecast", 46073, "OFF"),
  @(45, 46073.89583333334, 105.79, "forecast", 46073, "OFF"),
  @(46, 46073.91666666666, 105.79, "forecast", 46073, "OFF"),
  @(47, 46073.9375, 95.78986999999999, "forecast", 46073, "OFF"),
  @(48, 46073.95833333334, 95.32526, "forecast", 46073, "OFF"),
  @(49, 46073.97916666666, 105.2107, "forecast", 46073, "OFF"),
  @(50, 46074, 105.79, "forecast", 46074, "OFF"),
  @(51, 46074.02083333334, 97.26833000000001, "forecast", 46074, "OFF"),
  @(52, 46074.04166666666, 96.72154999999999, "forecast", 46074, "OFF"),
  @(53, 46074.0625, 100.00266, "forecast", 46074, "OFF"),
  @(54, 46074.08333333334, 102.09916, "forecast", 46074, "OFF"),
  @(55, 46074.10416666666, 97.38017000000001, "forecast", 46074, "OFF"),
  @(56, 46074.125, 97.2698, "forecast", 46074, "OFF"),
  @(57, 46074.14583333334, 96.72969000000001, "forecast", 46074, "OFF"),
  @(58, 46074.16666666666, 87.48866, "forecast", 46074, "OFF"),
  @(59, 46074.1875, 91.07102, "forecast", 46074, "OFF"),
  @(60, 46074.20833333334, 95.51734, "forecast", 46074, "OFF"),
  @(61, 46074.22916666666, 100.37296, "forecast", 46074, "OFF"),
  @(62, 46074.25, 105, "forecast", 46074, "OFF"),
  @(63, 46074.27083333334, 91.99872999999999, "forecast", 46074, "OFF"),
  @(64, 46074.29166666666, 75.02657000000001, "forecast", 46074, "ON"),
  @(65, 46074.3125, 19.62042, "forecast", 46074, "ON"),
  @(66, 46074.33333333334, 6.88598, "forecast", 46074, "ON"),
  @(67, 46074.35416666666, 0.7, "forecast", 46074, "ON"),
  @(68, 46074.375, 0.51, "forecast", 46074, "ON"),
  @(69, 46074.39583333334, 0.0003, "forecast", 46074, "ON"),
  @(70, 46074.41666666666, -3.11157, "forecast", 46074, "ON"),
  @(71, 46074.4375, 0.51, "forecast", 46074, "ON"),
  @(72, 46074.45833333334, 0.51, "forecast", 46074, "ON"),
  @(73, 46074.47916666666, 0.51, "forecast", 46074, "ON"),
  @(74, 46074.5, 0.7, "forecast", 46074, "ON"),
  @(75, 46074.52083333334, 0.7, "forecast", 46074, "ON"),
  @(76, 46074.54166666666, 0.51, "forecast", 46074, "ON"),
  @(77, 46074.5625, 11.70505, "forecast", 46074, "ON"),
  @(78, 46074.58333333334, 37.89, "forecast", 46074, "ON"),
  @(79, 46074.60416666666, 35.88, "forecast", 46074, "ON"),
  @(80, 46074.625, 37.89, "forecast", 46074, "ON"),
  @(81, 46074.64583333334, 37.89, "forecast", 46074, "ON"),
  @(82, 46074.66666666666, 37.89, "forecast", 46074, "ON"),
  @(83, 46074.6875, 37.89, "forecast", 46074, "ON"),
  @(84, 46074.70833333334, 56.75952, "forecast", 46074, "ON"),
  @(85, 46074.72916666666, 58.38328, "forecast", 46074, "ON"),
  @(86, 46074.75, 97.77654, "forecast", 46074, "OFF"),
  @(87, 46074.77083333334, 108.01, "forecast", 46074, "OFF"),
  @(88, 46074.79166666666, 128.30028, "forecast", 46074, "OFF"),
  @(89, 46074.8125, 108.84959, "forecast", 46074, "OFF"),
  @(90, 46074.83333333334, 108.01, "forecast", 46074, "OFF"),
  @(91, 46074.85416666666, 108.01, "forecast", 46074, "OFF"),
  @(92, 46074.875, 102.46052, "forecast", 46074, "OFF"),
  @(93, 46074.89583333334, 84.79000000000001, "forecast", 46074, "OFF"),
  @(94, 46074.91666666666, 78, "forecast", 46074, "OFF"),
  @(95, 46074.9375, 78, "forecast", 46074, "OFF"),
  @(96, 46074.95833333334, 84.79000000000001, "forecast", 46074, "OFF"),
  @(97, 46074.97916666666, 84.79000000000001, "forecast", 46074, "OFF")
)

foreach ($row in $detailedData) {
    $r = $row[0]
    $wsDetailed.Cells.Item($r, 1).Value = $row[1]
    $wsDetailed.Cells.Item($r, 1).NumberFormat = $fmtA_det
    $wsDetailed.Cells.Item($r, 2).Value = $row[2]
    $wsDetailed.Cells.Item($r, 3).Value = $row[3]
    $wsDetailed.Cells.Item($r, 4).Value = $row[4]
    $wsDetailed.Cells.Item($r, 4).NumberFormat = $fmtD_det
    $wsDetailed.Cells.Item($r, 5).Value = $row[5]
}
